$d = $word.ActiveDocument

# Locate the bold/underlined "Terminal" heading paragraph, then the
# (empty, bookmark-carrying) paragraph right after it -- that is the
# paragraph whose "sz 52 / szCs 52" run formatting must be dropped and
# whose bookmark must move down into the new body paragraph.
$headingRange = $d.Content
$found = $headingRange.Find.Execute("Terminal", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Start -le $headingRange.Start -and $candidate.Range.End -ge $headingRange.End) {
        $headingIndex = $i
        break
    }
}

$bookmarkParaIndex = $headingIndex + 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

$rangeStart = $bookmarkPara.Range.Start
$rangeEnd = $bookmarkPara.Range.End

# Drop the whole "sz 52 / szCs 52 / bookmark" paragraph -- it gets
# rebuilt (without the font size) below, and the bookmark is re-added
# to the new paragraph that holds the "Terminal only has..." text.
$deleteRange = $d.Range($rangeStart, $rangeEnd)
$deleteRange.Delete()

# Re-insert: an empty bold/underlined paragraph (no explicit size),
# a blank paragraph, and the new explanatory paragraph (the trailing
# paragraph of the payload merges into the existing "sz 28" paragraph
# that used to be empty, picking up the new runs + the _GoBack bookmark).
$insertPoint = $d.Range($rangeStart, $rangeStart)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p/><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Terminal only has a few things that weren’t able to be implemented in the final </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>ncurses</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> version. Using </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>“./</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">help” or “./commands” can bring up a list of commands. They all work, except for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>for</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>“./</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>changeroom”. The command is carried out and the client’s chatroom is changed, but it stops the server from routing messages.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml) | Out-Null
